# Insert a new record row at row 202 of the data table, which pushes all
# existing rows 202..304 down to 203..305 (Excel's native Insert behaviour
# copies the values/formatting of the shifted rows automatically), then
# populate the new row 202 with its own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(202).Insert()

$ws.Range("A202").Value = 9
$ws.Range("B202").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C202").Value = "Metropolitana"
$ws.Range("D202").Value = 44529
$ws.Range("E202").Value = 13
$ws.Range("F202").Value = 100112012
$ws.Range("G202").Value = "Espinaca"
$ws.Range("H202").Value = "Sin especificar"
$ws.Range("I202").Value = "Primera"
$ws.Range("J202").Value = 160
$ws.Range("K202").Value = 6000
$ws.Range("L202").Value = 7000
$ws.Range("M202").Value = 6500
$ws.Range("N202").Value = "$/cuna 10 kilos"
$ws.Range("O202").Value = "Provincia de Chacabuco"
$ws.Range("P202").Value = 650
$ws.Range("Q202").Value = 10
$ws.Range("R202").Value = "Hortaliza"
